$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two JST-connector line items (rows 3 and 4) are being dropped from
# the BOM; deleting the whole rows shifts the remaining rows (Honeywell
# sensor, capacitor, PCB, total) up by two, which also re-numbers the
# formulas (SUM range) and drops the now-unused shared strings.
$ws.Rows("3:4").Delete()

# Row deletion does not carry the worksheet's <hyperlinks> entries along
# with the shifted cells, so the two surviving links (Mouser datasheet on
# the Honeywell row, Digikey on the capacitor row) need to be reattached
# to their new K3/K4 home cells.
$ws.Range("K3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("K3"), "https://www.mouser.com/ProductDetail/Honeywell/ABPDANT030PG0D3?qs=%2Fha2pyFaduiufrMS3AAFabWwScpbL%252BOOUI4y8%252BgcuVXeHk%252B%2Fiw0hiw%3D%3D") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K4"), "https://www.digikey.com/product-detail/en/yageo/CC1206KKX7R0BB104/311-1488-1-ND/2833794") | Out-Null

# Re-adding a hyperlink resets the cell to Excel's built-in Hyperlink
# style; put it back to the workbook's existing hyperlink formatting.
$ws.Range("K3").Style = "Hyperlink"
$ws.Range("K4").Style = "Hyperlink"

# Match the cursor position left behind in the saved file.
$ws.Range("H12").Select()
